$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 44533c18-eeca-4e90-af39-4b8a6044c26d.md is now ready for handoff ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-24 02:33:53"

# --- zh-cn sheet: same file's status/handoff datetime updated ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-24 02:33:44"

# --- de-de sheet: same file's status/handoff datetime updated ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-24 02:33:53"
